$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, copying formatting (bold, border, alignment) from
# the neighboring header cell G1 ("sum"), then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill H2:H6 with 0 (numeric, no special style - matches data columns)
$ws.Range("H2:H6").Value = 0
